# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# Updates numeric columns H-N (price/profit calculations) across the ALC, ARM, BSM,
# CRP, CUL, GSM, LTW and WVR sheets, matching the upstream data refresh commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 453
$ws.Cells.Item(8, 9).Value = 58.88889
$ws.Cells.Item(8, 10).Value = 4000
$ws.Cells.Item(8, 11).Value = 176.66667
$ws.Cells.Item(8, 12).Value = 12000
$ws.Cells.Item(8, 13).Value = -37.66667000000001
$ws.Cells.Item(8, 14).Value = -12278

$ws.Cells.Item(15, 8).Value = 901.05554
$ws.Cells.Item(15, 9).Value = 901.05554
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 2703.16662
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = -2534.16662

$ws.Cells.Item(28, 8).Value = 590.25
$ws.Cells.Item(28, 9).Value = 581.5
$ws.Cells.Item(28, 10).Value = 625.25
$ws.Cells.Item(28, 11).Value = 581.5
$ws.Cells.Item(28, 12).Value = 625.25
$ws.Cells.Item(28, 13).Value = -96.5
$ws.Cells.Item(28, 14).Value = -1595.25

$ws.Cells.Item(70, 8).Value = 2265.9119
$ws.Cells.Item(70, 9).Value = 1797.0435
$ws.Cells.Item(70, 10).Value = 3246.2727
$ws.Cells.Item(70, 11).Value = 5391.1305
$ws.Cells.Item(70, 12).Value = 9738.8181
$ws.Cells.Item(70, 13).Value = -5121.1305
$ws.Cells.Item(70, 14).Value = -10278.8181

$ws.Cells.Item(73, 8).Value = 2265.9119
$ws.Cells.Item(73, 9).Value = 1797.0435
$ws.Cells.Item(73, 10).Value = 3246.2727
$ws.Cells.Item(73, 11).Value = 5391.1305
$ws.Cells.Item(73, 12).Value = 9738.8181
$ws.Cells.Item(73, 13).Value = -4455.1305
$ws.Cells.Item(73, 14).Value = -11610.8181

$ws.Cells.Item(93, 8).Value = 24467.742
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 24467.742
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 24467.742
$ws.Cells.Item(93, 14).Value = -29459.742

$ws.Cells.Item(106, 8).Value = 1990
$ws.Cells.Item(106, 9).Value = 1990
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 11).Value = 1990
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 13).Value = -1359
$ws.Cells.Item(106, 14).ClearContents()

$ws.Cells.Item(112, 8).Value = 1253.6984
$ws.Cells.Item(112, 9).Value = 350
$ws.Cells.Item(112, 10).Value = 1314.9661
$ws.Cells.Item(112, 11).Value = 1050
$ws.Cells.Item(112, 12).Value = 3944.8983
$ws.Cells.Item(112, 13).Value = 58
$ws.Cells.Item(112, 14).Value = -6160.898300000001

$ws.Cells.Item(115, 8).Value = 1329.1666
$ws.Cells.Item(115, 9).Value = 1177.2727
$ws.Cells.Item(115, 10).Value = 3000
$ws.Cells.Item(115, 11).Value = 3531.8181
$ws.Cells.Item(115, 12).Value = 9000
$ws.Cells.Item(115, 13).Value = -1964.8181

$ws.Cells.Item(128, 8).Value = 41835
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 41835
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 41835
$ws.Cells.Item(128, 14).Value = -51795

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(139, 8).Value = 36957.94
$ws.Cells.Item(139, 9).Value = 20000
$ws.Cells.Item(139, 10).Value = 38017.812
$ws.Cells.Item(139, 11).Value = 20000
$ws.Cells.Item(139, 12).Value = 38017.812
$ws.Cells.Item(139, 13).Value = -14860
$ws.Cells.Item(139, 14).Value = -48297.812

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1794.762
$ws.Cells.Item(107, 9).Value = 1233.9286
$ws.Cells.Item(107, 10).Value = 2916.4285
$ws.Cells.Item(107, 11).Value = 1233.9286
$ws.Cells.Item(107, 12).Value = 2916.4285
$ws.Cells.Item(107, 13).Value = 686.0714
$ws.Cells.Item(107, 14).Value = -6756.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 9775.571
$ws.Cells.Item(12, 9).Value = 7238.3335
$ws.Cells.Item(12, 10).Value = 24999
$ws.Cells.Item(12, 11).Value = 7238.3335
$ws.Cells.Item(12, 12).Value = 24999
$ws.Cells.Item(12, 13).Value = -7068.3335
$ws.Cells.Item(12, 14).Value = -25339

$ws.Cells.Item(31, 8).Value = 2638.6765
$ws.Cells.Item(31, 9).Value = 968.2727
$ws.Cells.Item(31, 10).Value = 5701.0835
$ws.Cells.Item(31, 11).Value = 968.2727
$ws.Cells.Item(31, 12).Value = 5701.0835
$ws.Cells.Item(31, 13).Value = -673.2727
$ws.Cells.Item(31, 14).Value = -6291.0835

$ws.Cells.Item(34, 8).Value = 2638.6765
$ws.Cells.Item(34, 9).Value = 968.2727
$ws.Cells.Item(34, 10).Value = 5701.0835
$ws.Cells.Item(34, 11).Value = 968.2727
$ws.Cells.Item(34, 12).Value = 5701.0835
$ws.Cells.Item(34, 13).Value = -766.2727
$ws.Cells.Item(34, 14).Value = -6105.0835

$ws.Cells.Item(58, 8).Value = 2874.9553
$ws.Cells.Item(58, 9).Value = 1713.8679
$ws.Cells.Item(58, 10).Value = 7270.5
$ws.Cells.Item(58, 11).Value = 1713.8679
$ws.Cells.Item(58, 12).Value = 7270.5
$ws.Cells.Item(58, 13).Value = -1510.8679
$ws.Cells.Item(58, 14).Value = -7676.5

$ws.Cells.Item(105, 8).Value = 1582.7778
$ws.Cells.Item(105, 9).Value = 1246.625
$ws.Cells.Item(105, 10).Value = 2071.7273
$ws.Cells.Item(105, 11).Value = 1246.625
$ws.Cells.Item(105, 12).Value = 2071.7273
$ws.Cells.Item(105, 13).Value = 500.375

$ws.Cells.Item(132, 8).Value = 3134.5652
$ws.Cells.Item(132, 9).Value = 1806.8
$ws.Cells.Item(132, 10).Value = 5624.125
$ws.Cells.Item(132, 11).Value = 5420.4
$ws.Cells.Item(132, 12).Value = 16872.375
$ws.Cells.Item(132, 13).Value = -2890.4
$ws.Cells.Item(132, 14).Value = -21932.375

$ws.Cells.Item(136, 8).Value = 2874.9553
$ws.Cells.Item(136, 9).Value = 1713.8679
$ws.Cells.Item(136, 10).Value = 7270.5
$ws.Cells.Item(136, 11).Value = 5141.6037
$ws.Cells.Item(136, 12).Value = 21811.5
$ws.Cells.Item(136, 13).Value = -2591.6037
$ws.Cells.Item(136, 14).Value = -26911.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 265.44446
$ws.Cells.Item(10, 9).Value = 265.44446
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 796.33338
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = -657.33338
$ws.Cells.Item(10, 14).ClearContents()

$ws.Cells.Item(75, 8).Value = 5878.5
$ws.Cells.Item(75, 9).Value = 500
$ws.Cells.Item(75, 10).Value = 7671.3335
$ws.Cells.Item(75, 11).Value = 1500
$ws.Cells.Item(75, 12).Value = 23014.0005
$ws.Cells.Item(75, 13).Value = -502
$ws.Cells.Item(75, 14).Value = -25010.0005

$ws.Cells.Item(78, 8).Value = 5878.5
$ws.Cells.Item(78, 9).Value = 500
$ws.Cells.Item(78, 10).Value = 7671.3335
$ws.Cells.Item(78, 11).Value = 4500
$ws.Cells.Item(78, 12).Value = 69042.0015
$ws.Cells.Item(78, 13).Value = 492
$ws.Cells.Item(78, 14).Value = -79026.0015

$ws.Cells.Item(122, 8).Value = 2925
$ws.Cells.Item(122, 9).Value = 854.8182
$ws.Cells.Item(122, 10).Value = 3835.88
$ws.Cells.Item(122, 11).Value = 7693.3638
$ws.Cells.Item(122, 12).Value = 34522.92
$ws.Cells.Item(122, 13).Value = -5243.3638
$ws.Cells.Item(122, 14).Value = -39422.92

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 11111805
$ws.Cells.Item(107, 9).Value = 271
$ws.Cells.Item(107, 10).Value = 18519494
$ws.Cells.Item(107, 11).Value = 271
$ws.Cells.Item(107, 12).Value = 18519494
$ws.Cells.Item(107, 13).Value = 1649
$ws.Cells.Item(107, 14).Value = -18523334

$ws.Cells.Item(126, 8).Value = 3219.01
$ws.Cells.Item(126, 9).Value = 2823.9487
$ws.Cells.Item(126, 10).Value = 4686.381
$ws.Cells.Item(126, 11).Value = 8471.846099999999
$ws.Cells.Item(126, 12).Value = 14059.143
$ws.Cells.Item(126, 13).Value = -6001.846099999999
$ws.Cells.Item(126, 14).Value = -18999.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(12, 8).Value = 300
$ws.Cells.Item(12, 9).Value = 300
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 300
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = -130
$ws.Cells.Item(12, 14).ClearContents()

$ws.Cells.Item(46, 8).Value = 1803.8572
$ws.Cells.Item(46, 9).Value = 1410.6316
$ws.Cells.Item(46, 10).Value = 2634
$ws.Cells.Item(46, 11).Value = 1410.6316
$ws.Cells.Item(46, 12).Value = 2634
$ws.Cells.Item(46, 13).Value = -1222.6316
$ws.Cells.Item(46, 14).Value = -3010

$ws.Cells.Item(61, 8).Value = 1869.9412
$ws.Cells.Item(61, 9).Value = 1791.1538
$ws.Cells.Item(61, 10).Value = 2126
$ws.Cells.Item(61, 11).Value = 1791.1538
$ws.Cells.Item(61, 12).Value = 2126
$ws.Cells.Item(61, 13).Value = -1589.1538

$ws.Cells.Item(74, 8).Value = 40843.75
$ws.Cells.Item(74, 9).Value = 24625
$ws.Cells.Item(74, 10).Value = 46250
$ws.Cells.Item(74, 11).Value = 24625
$ws.Cells.Item(74, 12).Value = 46250
$ws.Cells.Item(74, 13).Value = -23627

$ws.Cells.Item(77, 8).Value = 40843.75
$ws.Cells.Item(77, 9).Value = 24625
$ws.Cells.Item(77, 10).Value = 46250
$ws.Cells.Item(77, 11).Value = 73875
$ws.Cells.Item(77, 12).Value = 138750
$ws.Cells.Item(77, 13).Value = -68883

$ws.Cells.Item(113, 8).Value = 1869.9412
$ws.Cells.Item(113, 9).Value = 1791.1538
$ws.Cells.Item(113, 10).Value = 2126
$ws.Cells.Item(113, 11).Value = 1791.1538
$ws.Cells.Item(113, 12).Value = 2126
$ws.Cells.Item(113, 13).Value = 378.8462

$ws.Cells.Item(132, 8).Value = 3855.111
$ws.Cells.Item(132, 9).Value = 2572.7896
$ws.Cells.Item(132, 10).Value = 6900.625
$ws.Cells.Item(132, 11).Value = 7718.3688
$ws.Cells.Item(132, 12).Value = 20701.875
$ws.Cells.Item(132, 13).Value = -5188.3688
$ws.Cells.Item(132, 14).Value = -25761.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(19, 8).Value = 9990
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 9990
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 9990
$ws.Cells.Item(19, 14).Value = -10338

$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).ClearContents()

$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 11).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 6708.3125
$ws.Cells.Item(113, 9).Value = 9295.727999999999
$ws.Cells.Item(113, 10).Value = 1016
$ws.Cells.Item(113, 11).Value = 27887.184
$ws.Cells.Item(113, 12).Value = 3048
$ws.Cells.Item(113, 13).Value = -25717.184

$ws.Cells.Item(122, 8).Value = 4555.6665
$ws.Cells.Item(122, 9).Value = 2102
$ws.Cells.Item(122, 10).Value = 6308.2856
$ws.Cells.Item(122, 11).Value = 6306
$ws.Cells.Item(122, 12).Value = 18924.8568
$ws.Cells.Item(122, 13).Value = -3856
$ws.Cells.Item(122, 14).Value = -23824.8568

$ws.Cells.Item(126, 8).Value = 1067446.2
$ws.Cells.Item(126, 9).Value = 2652
$ws.Cells.Item(126, 10).Value = 1777309
$ws.Cells.Item(126, 11).Value = 7956
$ws.Cells.Item(126, 12).Value = 5331927
$ws.Cells.Item(126, 13).Value = -5486
$ws.Cells.Item(126, 14).Value = -5336867

$ws.Cells.Item(127, 8).Value = 39930
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 39930
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 39930
$ws.Cells.Item(127, 14).Value = -49850

$ws.Cells.Item(128, 8).Value = 41812.145
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 41812.145
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 41812.145
$ws.Cells.Item(128, 14).Value = -51772.145

$ws.Cells.Item(136, 8).Value = 11385
$ws.Cells.Item(136, 9).Value = 11296.7
$ws.Cells.Item(136, 10).Value = 11511.143
$ws.Cells.Item(136, 11).Value = 33890.10000000001
$ws.Cells.Item(136, 12).Value = 34533.429
$ws.Cells.Item(136, 13).Value = -31340.10000000001

$ws.Cells.Item(137, 8).Value = 49257.5
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 49257.5
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 49257.5
$ws.Cells.Item(137, 14).Value = -59457.5

$ws.Cells.Item(139, 8).Value = 36912.777
$ws.Cells.Item(139, 9).Value = 40650
$ws.Cells.Item(139, 10).Value = 36445.625
$ws.Cells.Item(139, 11).Value = 40650
$ws.Cells.Item(139, 12).Value = 36445.625
$ws.Cells.Item(139, 13).Value = -35510
$ws.Cells.Item(139, 14).Value = -46725.625
